$wb = $excel.ActiveWorkbook

# --- ALC row 116 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3483.8857
$ws.Range("I116").Value = 3940.8096
$ws.Range("J116").Value = 2798.5
$ws.Range("K116").Value = 3940.8096
$ws.Range("L116").Value = 2798.5
$ws.Range("M116").Value = -498.8096
$ws.Range("N116").Value = -9682.5

# --- ALC row 132 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2916.3333
$ws.Range("I132").Value = 1321.2916
$ws.Range("J132").Value = 6562.143
$ws.Range("K132").Value = 3963.8748
$ws.Range("L132").Value = 19686.429
$ws.Range("M132").Value = -1433.8748
$ws.Range("N132").Value = -24746.429

# --- ALC row 138 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1519.9375
$ws.Range("I138").Value = 980.86487
$ws.Range("J138").Value = 2258.6667
$ws.Range("K138").Value = 2942.59461
$ws.Range("L138").Value = 6776.000100000001
$ws.Range("M138").Value = 2197.40539
$ws.Range("N138").Value = -17056.0001

# --- ARM row 2 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2106.0476
$ws.Range("I2").Value = 2316.2144
$ws.Range("J2").Value = 1685.7142
$ws.Range("K2").Value = 2316.2144
$ws.Range("L2").Value = 1685.7142
$ws.Range("M2").Value = -2203.2144
$ws.Range("N2").Value = -1911.7142

# --- ARM row 32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 929.74
$ws.Range("I32").Value = 886.18555
$ws.Range("J32").Value = 2338
$ws.Range("K32").Value = 886.18555
$ws.Range("L32").Value = 2338
$ws.Range("M32").Value = -599.18555
$ws.Range("N32").Value = -2912

# --- ARM row 61 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 422742.88
$ws.Range("I61").Value = 325326.53
$ws.Range("K61").Value = 325326.53
$ws.Range("M61").Value = -325114.53

# --- ARM row 116 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2106.0476
$ws.Range("I116").Value = 2316.2144
$ws.Range("J116").Value = 1685.7142
$ws.Range("K116").Value = 2316.2144
$ws.Range("L116").Value = 1685.7142
$ws.Range("M116").Value = -22.21439999999984
$ws.Range("N116").Value = -6273.7142

# --- ARM row 122 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6404.952
$ws.Range("I122").Value = 6533.5557
$ws.Range("J122").Value = 5633.3335
$ws.Range("K122").Value = 19600.6671
$ws.Range("L122").Value = 16900.0005
$ws.Range("M122").Value = -17150.6671
$ws.Range("N122").Value = -21800.0005

# --- ARM row 132 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3022.5686
$ws.Range("I132").Value = 2843.9167
$ws.Range("J132").Value = 3451.3333
$ws.Range("K132").Value = 8531.750100000001
$ws.Range("L132").Value = 10353.9999
$ws.Range("M132").Value = -6001.750100000001
$ws.Range("N132").Value = -15413.9999

# --- ARM row 136 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 422742.88
$ws.Range("I136").Value = 325326.53
$ws.Range("K136").Value = 975979.5900000001
$ws.Range("M136").Value = -973429.5900000001

# --- BSM row 3 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2106.0476
$ws.Range("I3").Value = 2316.2144
$ws.Range("J3").Value = 1685.7142
$ws.Range("K3").Value = 2316.2144
$ws.Range("L3").Value = 1685.7142
$ws.Range("M3").Value = -2202.2144
$ws.Range("N3").Value = -1913.7142

# --- BSM row 22 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 250
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -846

# --- BSM row 99 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7974421
$ws.Range("I99").Value = 2569225.5
$ws.Range("J99").Value = 35000400
$ws.Range("K99").Value = 2569225.5
$ws.Range("L99").Value = 35000400
$ws.Range("M99").Value = -2567727.5
$ws.Range("N99").Value = -35003396

# --- BSM row 133 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 39995
$ws.Range("J133").Value = 39995
$ws.Range("L133").Value = 39995
$ws.Range("N133").Value = -50115

# --- CRP row 31 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2557.5186
$ws.Range("I31").Value = 1670.591
$ws.Range("K31").Value = 1670.591
$ws.Range("M31").Value = -1375.591

# --- CRP row 34 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2557.5186
$ws.Range("I34").Value = 1670.591
# J34 unchanged (6460 -> 6460), skipping
$ws.Range("K34").Value = 1670.591
$ws.Range("M34").Value = -1468.591

# --- CRP row 58 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2743.9
$ws.Range("I58").Value = 2910
$ws.Range("J58").Value = 2143.3845
$ws.Range("K58").Value = 2910
$ws.Range("L58").Value = 2143.3845
$ws.Range("M58").Value = -2707
$ws.Range("N58").Value = -2549.3845

# --- CRP row 132 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2425.2173
$ws.Range("I132").Value = 1341.4286
$ws.Range("J132").Value = 4111.1113
$ws.Range("K132").Value = 4024.2858
$ws.Range("L132").Value = 12333.3339
$ws.Range("M132").Value = -1494.2858
$ws.Range("N132").Value = -17393.3339

# --- CRP row 133 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 40975.2
$ws.Range("I133").Value = 21296
$ws.Range("J133").Value = 45895
$ws.Range("K133").Value = 21296
$ws.Range("L133").Value = 45895
$ws.Range("M133").Value = -18766
$ws.Range("N133").Value = -50955

# --- CRP row 136 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2743.9
$ws.Range("I136").Value = 2910
$ws.Range("J136").Value = 2143.3845
$ws.Range("K136").Value = 8730
$ws.Range("L136").Value = 6430.1535
$ws.Range("M136").Value = -6180
$ws.Range("N136").Value = -11530.1535

# --- CUL row 14 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 69
$ws.Range("I14").Value = 69
$ws.Range("K14").Value = 207
$ws.Range("M14").Value = -34

# --- CUL row 107 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 501.95
$ws.Range("I107").Value = 331.75
$ws.Range("J107").Value = 615.4167
$ws.Range("K107").Value = 995.25
$ws.Range("L107").Value = 1846.2501
$ws.Range("M107").Value = 924.75
$ws.Range("N107").Value = -5686.2501

# --- GSM row 132 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4050.6365
$ws.Range("I132").Value = 4447.579
$ws.Range("J132").Value = 3511.9285
$ws.Range("K132").Value = 13342.737
$ws.Range("L132").Value = 10535.7855
$ws.Range("M132").Value = -10812.737
$ws.Range("N132").Value = -15595.7855

# --- LTW row 22 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 916.4545000000001
$ws.Range("I22").Value = 947.2857
$ws.Range("K22").Value = 947.2857
$ws.Range("M22").Value = -652.2857

# --- LTW row 27 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 916.4545000000001
$ws.Range("I27").Value = 947.2857
$ws.Range("K27").Value = 947.2857
$ws.Range("M27").Value = -840.2857

# --- LTW row 55 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 302.05
$ws.Range("I55").Value = 177.5625
$ws.Range("J55").Value = 800
$ws.Range("K55").Value = 177.5625
$ws.Range("L55").Value = 800
$ws.Range("M55").Value = -4.5625
$ws.Range("N55").Value = -1146

# --- LTW row 132 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6420.38
$ws.Range("I132").Value = 2125.2974
$ws.Range("K132").Value = 6375.8922
$ws.Range("M132").Value = -3845.8922

# --- WVR row 8 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 4900
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 4900
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 4900
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -5180

# --- WVR row 101 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 14534
$ws.Range("J101").Value = 14534
$ws.Range("L101").Value = 14534
$ws.Range("N101").Value = -21024

# --- WVR row 132 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1674.6327
$ws.Range("I132").Value = 1208.0667
$ws.Range("K132").Value = 3624.2001
$ws.Range("M132").Value = -1094.2001
